# Apply data update to the "Inscricoes" worksheet, matching the commit
# "Data update using git" diff: update Pagos (F), Inscrições homologadas (H)
# and Inscritos (E) figures for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 3
$ws.Range("F3").Value = 14
$ws.Range("H3").Value = 17

# Row 7
$ws.Range("F7").Value = 20
$ws.Range("H7").Value = 21

# Row 8
$ws.Range("E8").Value = 45

# Row 13
$ws.Range("F13").Value = 6
$ws.Range("H13").Value = 6

# Row 14
$ws.Range("E14").Value = 37

# Row 16
$ws.Range("F16").Value = 89
$ws.Range("H16").Value = 176

# Row 17
$ws.Range("E17").Value = 25

$wb.Save()
